$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44644
$ws.Cells.Item(2, 10).Value = 140
$ws.Cells.Item(2, 11).Value = 2500
$ws.Cells.Item(2, 12).Value = 3000
$ws.Cells.Item(2, 13).Value = 2786
$ws.Cells.Item(2, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(2, 16).Value = 464

$ws.Cells.Item(3, 4).Value = 44957
$ws.Cells.Item(3, 10).Value = 70
$ws.Cells.Item(3, 11).Value = 1500
$ws.Cells.Item(3, 12).Value = 2000
$ws.Cells.Item(3, 13).Value = 1857
$ws.Cells.Item(3, 15).Value = "Región Metropolitana"
$ws.Cells.Item(3, 16).Value = 310

$ws.Cells.Item(4, 4).Value = 44630
$ws.Cells.Item(4, 10).Value = 90
$ws.Cells.Item(4, 11).Value = 2500
$ws.Cells.Item(4, 12).Value = 3000
$ws.Cells.Item(4, 13).Value = 2722
$ws.Cells.Item(4, 15).Value = "Región Metropolitana"
$ws.Cells.Item(4, 16).Value = 454

$ws.Cells.Item(5, 4).Value = 44659
$ws.Cells.Item(5, 10).Value = 90
$ws.Cells.Item(5, 11).Value = 2500
$ws.Cells.Item(5, 12).Value = 3000
$ws.Cells.Item(5, 13).Value = 2722
$ws.Cells.Item(5, 15).Value = "Región Metropolitana"
$ws.Cells.Item(5, 16).Value = 454

$ws.Cells.Item(6, 4).Value = 44658
$ws.Cells.Item(6, 10).Value = 180
$ws.Cells.Item(6, 11).Value = 2500
$ws.Cells.Item(6, 12).Value = 3000
$ws.Cells.Item(6, 13).Value = 2778
$ws.Cells.Item(6, 15).Value = "Región Metropolitana"
$ws.Cells.Item(6, 16).Value = 463

$ws.Cells.Item(7, 4).Value = 44876
$ws.Cells.Item(7, 10).Value = 80
$ws.Cells.Item(7, 11).Value = 6500
$ws.Cells.Item(7, 12).Value = 7000
$ws.Cells.Item(7, 13).Value = 6812
$ws.Cells.Item(7, 15).Value = "Región Metropolitana"
$ws.Cells.Item(7, 16).Value = 1135

$ws.Cells.Item(8, 4).Value = 44643
$ws.Cells.Item(8, 10).Value = 90
$ws.Cells.Item(8, 11).Value = 2800
$ws.Cells.Item(8, 12).Value = 3000
$ws.Cells.Item(8, 13).Value = 2911
$ws.Cells.Item(8, 15).Value = "Región Metropolitana"
$ws.Cells.Item(8, 16).Value = 485

$ws.Cells.Item(9, 4).Value = 44672
$ws.Cells.Item(9, 10).Value = 140
$ws.Cells.Item(9, 11).Value = 3000
$ws.Cells.Item(9, 12).Value = 3500
$ws.Cells.Item(9, 13).Value = 3286
$ws.Cells.Item(9, 15).Value = "Región Metropolitana"
$ws.Cells.Item(9, 16).Value = 548

$ws.Cells.Item(10, 4).Value = 44650
$ws.Cells.Item(10, 10).Value = 130
$ws.Cells.Item(10, 11).Value = 3000
$ws.Cells.Item(10, 12).Value = 3500
$ws.Cells.Item(10, 13).Value = 3308
$ws.Cells.Item(10, 15).Value = "Región Metropolitana"
$ws.Cells.Item(10, 16).Value = 551

$ws.Cells.Item(11, 4).Value = 44671
$ws.Cells.Item(11, 10).Value = 150
$ws.Cells.Item(11, 11).Value = 3500
$ws.Cells.Item(11, 12).Value = 4000
$ws.Cells.Item(11, 13).Value = 3733
$ws.Cells.Item(11, 15).Value = "Región Metropolitana"
$ws.Cells.Item(11, 16).Value = 622

$ws.Cells.Item(12, 4).Value = 44685
$ws.Cells.Item(12, 10).Value = 150
$ws.Cells.Item(12, 11).Value = 3000
$ws.Cells.Item(12, 12).Value = 3500
$ws.Cells.Item(12, 13).Value = 3267
$ws.Cells.Item(12, 15).Value = "Región Metropolitana"
$ws.Cells.Item(12, 16).Value = 544

$ws.Cells.Item(13, 4).Value = 44631
$ws.Cells.Item(13, 10).Value = 110
$ws.Cells.Item(13, 11).Value = 3000
$ws.Cells.Item(13, 12).Value = 3500
$ws.Cells.Item(13, 13).Value = 3273
$ws.Cells.Item(13, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(13, 16).Value = 546

$ws.Cells.Item(14, 4).Value = 44637
$ws.Cells.Item(14, 10).Value = 170
$ws.Cells.Item(14, 11).Value = 2800
$ws.Cells.Item(14, 12).Value = 3000
$ws.Cells.Item(14, 13).Value = 2906
$ws.Cells.Item(14, 15).Value = "Región Metropolitana"
$ws.Cells.Item(14, 16).Value = 484

$ws.Cells.Item(15, 4).Value = 44987
$ws.Cells.Item(15, 10).Value = 130
$ws.Cells.Item(15, 11).Value = 4500
$ws.Cells.Item(15, 12).Value = 5000
$ws.Cells.Item(15, 13).Value = 4692
$ws.Cells.Item(15, 15).Value = "Región Metropolitana"
$ws.Cells.Item(15, 16).Value = 782

$ws.Cells.Item(16, 4).Value = 45021
$ws.Cells.Item(16, 10).Value = 50
$ws.Cells.Item(16, 11).Value = 4500
$ws.Cells.Item(16, 12).Value = 5000
$ws.Cells.Item(16, 13).Value = 4700
$ws.Cells.Item(16, 15).Value = "Región Metropolitana"
$ws.Cells.Item(16, 16).Value = 783
